# Daily Scores update - 2025-02-13
# Source workbook: data/Firmen/ipt/2025-02-01/Daily_Scores.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows 46-49 (2025-02-12 Sportfinke-league recompute) ---
# Row 46 (abs_activity, 2025-02-12)
$ws.Range("C46").Value = 5.123940646266563
$ws.Range("D46").Value = 6.191675377571315
$ws.Range("E46").Value = 6.831376101351728
$ws.Range("F46").Value = 9.642405490877502
$ws.Range("G46").Value = 9.470279925674344
$ws.Range("I46").Value = 9.86105701722066
$ws.Range("J46").Value = 8.984847155626849
$ws.Range("K46").Value = 8.881755699782785
$ws.Range("M46").Value = 9.11306104327611
$ws.Range("P46").Value = 49.28147043357219
$ws.Range("Q46").Value = 44.81892802407566

# Row 47 (rel_activity, 2025-02-12)
$ws.Range("F47").Value = 6.388162183790014
$ws.Range("K47").Value = 8.919100964265457
$ws.Range("L47").Value = 10
$ws.Range("M47").Value = 5.517927459966756
$ws.Range("P47").Value = 19.43702842423221
$ws.Range("Q47").Value = 45.63816218379002

# Row 48 (abs_sleep, 2025-02-12)
$ws.Range("K48").Value = 6.866666666666667
$ws.Range("P48").Value = 53.6

# Row 49 (rel_sleep, 2025-02-12)
$ws.Range("K49").Value = 0
$ws.Range("P49").Value = 19.99232158988257

# --- New rows 50-53 for 2025-02-13 ---
# Row 50 (abs_activity, 2025-02-13)
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "2025-02-13"
$ws.Range("B50").Value = "abs_activity"
$ws.Range("C50").Value = 6.773942304219698
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 8.689517589037555
$ws.Range("F50").Value = 9.19160777329574
$ws.Range("G50").Value = 9.034852341172611
$ws.Range("H50").Value = 10
$ws.Range("I50").Value = 10
$ws.Range("J50").Value = 4.063173608073821
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 10
$ws.Range("M50").Value = 5.148748490536102
$ws.Range("N50").Value = 0
$ws.Range("O50").Value = 0
$ws.Range("P50").Value = 39.64706072496596
$ws.Range("Q50").Value = 33.25478138136956

# Row 51 (rel_activity, 2025-02-13)
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = "2025-02-13"
$ws.Range("B51").Value = "rel_activity"
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 5
$ws.Range("E51").Value = 5.401848503114326
$ws.Range("F51").Value = 5.384024577572965
$ws.Range("G51").Value = 6.117991056071551
$ws.Range("H51").Value = 10
$ws.Range("I51").Value = 6.531583773439022
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 8.919100964265457
$ws.Range("L51").Value = 10
$ws.Range("M51").Value = 0
$ws.Range("N51").Value = 5
$ws.Range("O51").Value = 5
$ws.Range("P51").Value = 31.97052429689035
$ws.Range("Q51").Value = 35.38402457757297

# Row 52 (abs_sleep, 2025-02-13)
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = "2025-02-13"
$ws.Range("B52").Value = "abs_sleep"
$ws.Range("C52").Value = 8.066666666666666
$ws.Range("D52").Value = 10
$ws.Range("E52").Value = 10
$ws.Range("F52").Value = 8.4
$ws.Range("G52").Value = 9.4
$ws.Range("H52").Value = 3.733333333333334
$ws.Range("I52").Value = 10
$ws.Range("J52").Value = 10
$ws.Range("K52").Value = 10
$ws.Range("L52").Value = 4.666666666666668
$ws.Range("M52").Value = 0
$ws.Range("N52").Value = 0
$ws.Range("O52").Value = 0
$ws.Range("P52").Value = 47.46666666666667
$ws.Range("Q52").Value = 36.8

# Row 53 (rel_sleep, 2025-02-13)
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = "2025-02-13"
$ws.Range("B53").Value = "rel_sleep"
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 10
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 7.341147392396427
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 8.816216216216217
$ws.Range("K53").Value = 10
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("N53").Value = 0
$ws.Range("O53").Value = 0
$ws.Range("P53").Value = 10
$ws.Range("Q53").Value = 26.15736360861264

